# FIM run and forecast pull forward
# Update the "projection" rows (both "current" and "difference" sources)
# for the forecast quarters in columns P:X (2024 Q2 - 2026 Q2) with the
# refreshed FIM model output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P3").Value = 0.0085
$ws.Range("Q3").Value = 0.0071
$ws.Range("R3").Value = 0.0062

$ws.Range("P9").Value = -0.4557
$ws.Range("Q9").Value = -0.1629
$ws.Range("R9").Value = -0.1465

$ws.Range("P11").Value = -0.0256
$ws.Range("Q11").Value = -0.0302
$ws.Range("R11").Value = -0.0229
$ws.Range("S11").Value = -0.0161
$ws.Range("T11").Value = 0.0049
$ws.Range("U11").Value = -0.0031
$ws.Range("V11").Value = -0.0127
$ws.Range("W11").Value = -0.0274
$ws.Range("X11").Value = -1.7933

$ws.Range("P13").Value = -0.1279
$ws.Range("Q13").Value = 0.0584
$ws.Range("R13").Value = 0.0104
$ws.Range("S13").Value = 0.0486
$ws.Range("T13").Value = -0.0036
$ws.Range("U13").Value = 0.0157
$ws.Range("V13").Value = 0.033
$ws.Range("W13").Value = 0.0476
$ws.Range("X13").Value = -5.2081

$ws.Range("P15").Value = 0.1788
$ws.Range("Q15").Value = 0.3287
$ws.Range("R15").Value = 0.2853
$ws.Range("S15").Value = -0.0597
$ws.Range("T15").Value = -0.0771
$ws.Range("U15").Value = -0.046
$ws.Range("V15").Value = -0.0541
$ws.Range("W15").Value = 0.0282
$ws.Range("X15").Value = 7.6011

$ws.Range("T17").Value = -0.014

$ws.Range("T19").Value = -0.01
$ws.Range("V19").Value = -0.0077
$ws.Range("W19").Value = -0.0071

$ws.Range("P21").Value = 0.0095
$ws.Range("Q21").Value = 0.0082
$ws.Range("R21").Value = 0.0077

$ws.Range("P23").Value = 0.0152
$ws.Range("Q23").Value = 0.0421
$ws.Range("R23").Value = 0.0334
$ws.Range("S23").Value = -0.1148
$ws.Range("T23").Value = -0.0988
$ws.Range("U23").Value = -0.0906

$ws.Range("T25").Value = -0.0236

$ws.Range("P27").Value = -0.2931
$ws.Range("Q27").Value = -0.2649
$ws.Range("R27").Value = -0.1078
$ws.Range("S27").Value = -0.0597
$ws.Range("T27").Value = -0.0948
$ws.Range("U27").Value = -0.0707
$ws.Range("V27").Value = -0.0322
$ws.Range("W27").Value = -0.0131
$ws.Range("X27").Value = -0.0742

$ws.Range("P31").Value = -0.5629
$ws.Range("Q31").Value = -0.2064
$ws.Range("R31").Value = 0.094
$ws.Range("S31").Value = -0.4702
$ws.Range("T31").Value = -0.6931
$ws.Range("U31").Value = -0.4406
$ws.Range("V31").Value = -0.3915
$ws.Range("W31").Value = -0.282
$ws.Range("X31").Value = -73.7292

$ws.Range("P35").Value = -0.2316
$ws.Range("Q35").Value = -0.0305
$ws.Range("R35").Value = 0.0015

$ws.Range("P43").Value = 0.0493
$ws.Range("Q43").Value = -0.1489
$ws.Range("R43").Value = -0.1368

$ws.Range("P45").Value = -0.0165
$ws.Range("Q45").Value = -0.0088
$ws.Range("R45").Value = 0.0082
$ws.Range("S45").Value = 0.0227
$ws.Range("T45").Value = 0.0023
$ws.Range("U45").Value = 0.0025
$ws.Range("V45").Value = 0.0078
$ws.Range("W45").Value = 0.0127
$ws.Range("X45").Value = 0.0755

$ws.Range("P47").Value = 0.1152
$ws.Range("Q47").Value = 0.0565
$ws.Range("R47").Value = 0.1805
$ws.Range("S47").Value = 0.0073
$ws.Range("T47").Value = -0.0236
$ws.Range("U47").Value = -0.0165
$ws.Range("V47").Value = -0.0087
$ws.Range("X47").Value = -1.0114

$ws.Range("P49").Value = 0.2563
$ws.Range("Q49").Value = 0.175
$ws.Range("R49").Value = 0.1406
$ws.Range("S49").Value = 0.0987
$ws.Range("T49").Value = 0.0243
$ws.Range("U49").Value = 0.0542
$ws.Range("V49").Value = 0.0646
$ws.Range("W49").Value = 0.0267
$ws.Range("X49").Value = 3.9236

$ws.Range("P51").Value = 0.0088
$ws.Range("Q51").Value = 0.008
$ws.Range("R51").Value = 0.007

$ws.Range("P57").Value = -0.0014
$ws.Range("Q57").Value = -0.0003
$ws.Range("R57").Value = -0.0025
$ws.Range("S57").Value = -0.0034
$ws.Range("T57").Value = -0.0031
$ws.Range("U57").Value = -0.0024
$ws.Range("V57").Value = -0.0033
$ws.Range("W57").Value = -0.0031
$ws.Range("X57").Value = -0.0618

$ws.Range("P59").Value = 0.0017
$ws.Range("Q59").Value = 0.0011
$ws.Range("R59").Value = 0.0007

$ws.Range("P65").Value = -0.0959
$ws.Range("Q65").Value = -0.0956
$ws.Range("R65").Value = -0.0627

$ws.Range("P67").Value = 0.0947
$ws.Range("Q67").Value = 0.2079
$ws.Range("R67").Value = 0.2006
$ws.Range("S67").Value = 0.2678
$ws.Range("T67").Value = 0.2243
$ws.Range("U67").Value = 0.1675
$ws.Range("V67").Value = 0.0227
$ws.Range("W67").Value = 0.0867
$ws.Range("X67").Value = -1.366

$ws.Range("P69").Value = 0.0553
$ws.Range("Q69").Value = 0.0191
$ws.Range("R69").Value = 0.0115
$ws.Range("S69").Value = -0.0065
$ws.Range("T69").Value = -0.0484
$ws.Range("U69").Value = -0.0164
$ws.Range("V69").Value = -0.0155
$ws.Range("W69").Value = -0.0444
$ws.Range("X69").Value = -0.1839

$ws.Range("P71").Value = 0.0688
$ws.Range("Q71").Value = 0.0558
$ws.Range("R71").Value = 0.0767
$ws.Range("S71").Value = 0.0741
$ws.Range("T71").Value = 0.0746
$ws.Range("U71").Value = 0.0781
$ws.Range("V71").Value = 0.0633
$ws.Range("W71").Value = 0.0956
$ws.Range("X71").Value = 0.1382

$ws.Range("W73").Value = -0.0003

$ws.Range("Q75").Value = -0.0003
$ws.Range("R75").Value = -0.0003
$ws.Range("S75").Value = -0.0002
$ws.Range("W75").Value = -0.0002

$ws.Range("P77").Value = 0.0017
$ws.Range("Q77").Value = 0.0016
$ws.Range("R77").Value = 0.0007

$ws.Range("P79").Value = -0.0091
$ws.Range("Q79").Value = -0.0241
$ws.Range("R79").Value = -0.0365
$ws.Range("S79").Value = -0.0549
$ws.Range("T79").Value = -0.0441
$ws.Range("U79").Value = -0.0366

$ws.Range("R81").Value = -0.0003
$ws.Range("S81").Value = -0.0004

$ws.Range("P83").Value = -0.0011
$ws.Range("Q83").Value = -0.0021
$ws.Range("R83").Value = -0.0017
$ws.Range("S83").Value = -0.0015
$ws.Range("T83").Value = -0.0021
$ws.Range("U83").Value = -0.0019
$ws.Range("V83").Value = -0.0015
$ws.Range("W83").Value = -0.0015
$ws.Range("X83").Value = -0.0018

$ws.Range("P87").Value = 0.1109
$ws.Range("Q87").Value = 0.1567
$ws.Range("R87").Value = 0.213
$ws.Range("S87").Value = 0.2161
$ws.Range("T87").Value = 0.1264
$ws.Range("U87").Value = 0.1236
$ws.Range("V87").Value = -0.021
$ws.Range("W87").Value = -0.1911
$ws.Range("X87").Value = -3.082

$ws.Range("P91").Value = -0.021
$ws.Range("Q91").Value = -0.0198
$ws.Range("R91").Value = -0.0137

$ws.Range("P99").Value = -0.1065
$ws.Range("Q99").Value = -0.1052
$ws.Range("R99").Value = -0.0743

$ws.Range("P101").Value = -0.0004
$ws.Range("Q101").Value = -0.0003
$ws.Range("R101").Value = 0
$ws.Range("S101").Value = 0.0002
$ws.Range("T101").Value = 0.0001
$ws.Range("U101").Value = 0.0002
$ws.Range("V101").Value = 0.0004
$ws.Range("W101").Value = 0.0009
$ws.Range("X101").Value = 0.0026

$ws.Range("P103").Value = 0.0876
$ws.Range("Q103").Value = 0.084
$ws.Range("R103").Value = 0.0833
$ws.Range("S103").Value = 0.0133
$ws.Range("T103").Value = -0.0038
$ws.Range("U103").Value = -0.0012
$ws.Range("V103").Value = -0.0011
$ws.Range("X103").Value = -0.1176

$ws.Range("P105").Value = 0.0251
$ws.Range("Q105").Value = 0.0277
$ws.Range("R105").Value = 0.0298
$ws.Range("S105").Value = 0.0297
$ws.Range("T105").Value = 0.0296
$ws.Range("U105").Value = 0.032
$ws.Range("V105").Value = 0.0339
$ws.Range("W105").Value = 0.0487
$ws.Range("X105").Value = 0.092

$ws.Range("P107").Value = 0.0017
$ws.Range("Q107").Value = 0.0016
$ws.Range("R107").Value = 0.0007

$ws.Range("P113").Value = -0.0058
$ws.Range("Q113").Value = -0.0063
$ws.Range("R113").Value = -0.0082
$ws.Range("S113").Value = -0.0096
$ws.Range("T113").Value = -0.0094
$ws.Range("U113").Value = -0.0073
$ws.Range("V113").Value = -0.0057
$ws.Range("W113").Value = -0.0044
$ws.Range("X113").Value = 0.0094
